# Flights.xlsx — "Aircraft functionality added, and the import of the
# importFlights added."
#
# The diff adds five new flight rows (rows 6-10) below the existing four
# (rows 2-5), re-using several of the original shared strings (airline,
# plane, source/destination/dates that repeat) and introducing a handful
# of brand-new ones (Avianca, Taca, United, Delta, RF52, San Salvador,
# Lima Peru, new dates/times). It also re-confirms A4/A5 as plain
# row-index numbers 3/4, and leaves the final selection on the newly
# imported first row (A6:H6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-less inline pattern: the "FECHA" column holds plain text dates
# (not real Excel dates) in the source data, e.g. "6/10/21", "05/02/2021".
# Excel's smart-typing would otherwise silently convert those strings into
# date serial numbers, so each date cell is force-formatted as Text first
# and the format is cleared again afterwards (so no stray number format is
# left behind on the cell, matching the unstyled cells used elsewhere in
# this sheet).

# --- re-affirm the row index numbers for the pre-existing rows ----------
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# --- row 6 : imported flight #1 (re-uses nave31/rf56/san salvador/hawaii) -
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Avianca"
$ws.Range("C6").Value = "rf56"
$ws.Range("D6").Value = "san salvador"
$ws.Range("E6").Value = "hawaii"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "6/10/21"
$ws.Range("F6").ClearFormats()
$ws.Range("G6").Value = "2:00"
$ws.Range("H6").Value = "23:00"

# --- row 7 : imported flight #2 (re-uses nave5/hh7/Paris/new york) -------
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Taca"
$ws.Range("C7").Value = "hh7"
$ws.Range("D7").Value = "Paris"
$ws.Range("E7").Value = "new york"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "6/12/21"
$ws.Range("F7").ClearFormats()
$ws.Range("G7").Value = "7:00"
$ws.Range("H7").Value = "19:00"

# --- row 8 : imported flight #3 (re-uses SFDS3/Toronto/California/...) --
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "United"
$ws.Range("C8").Value = "SFDS3"
$ws.Range("D8").Value = "Toronto"
$ws.Range("E8").Value = "California"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "05/02/2021"
$ws.Range("F8").ClearFormats()
$ws.Range("G8").Value = "10:00"
$ws.Range("H8").Value = "2:00"

# --- row 9 : imported flight #4 (re-uses SDFD4/milan/sevillaa/...) ------
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Delta"
$ws.Range("C9").Value = "SDFD4"
$ws.Range("D9").Value = "milan"
$ws.Range("E9").Value = "sevillaa"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "12/05/2021"
$ws.Range("F9").ClearFormats()
$ws.Range("G9").Value = "15:00"
$ws.Range("H9").Value = "23:00"

# --- row 10 : imported flight #5 (brand-new aircraft/route) -------------
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Taca"
$ws.Range("C10").Value = "RF52"
$ws.Range("D10").Value = "San Salvador"
$ws.Range("E10").Value = "Lima, Peru"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "14/06/2021"
$ws.Range("F10").ClearFormats()
$ws.Range("G10").Value = "5:50"
$ws.Range("H10").Value = "10:00"

# --- leave the selection on the first newly-imported row -----------------
$ws.Range("A6:H6").Select()
